# Fix errors in the countries list:
# Insert a new row for "The Democratic Republic of the Congo" right before
# the existing "Democratic Republic of the Congo" row (row 323), pushing the
# later "Congo" name-variant rows (and everything below) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 323 (shifts 323..329 down to 324..330).
# Excel copies the formatting of the row above (row 322, style "2"),
# which matches the style used by the other Congo/DRC duplicate rows.
$ws.Rows(323).Insert()

# Fill in the new row with the same Iso3/Continent/SubContinent as the
# "Democratic Republic of the Congo" entry that used to sit at row 323.
$ws.Range("A323").Value = "The Democratic Republic of the Congo"
$ws.Range("B323").Value = "COD"
$ws.Range("C323").Value = "Africa"
$ws.Range("D323").Value = "Middle Africa"

# Restore the cursor/selection position recorded in the saved file.
$ws.Range("F322").Select() | Out-Null
